# Auto-generated Excel COM-interop script
# Updates cached market-price columns (H-N) on several Leve rows
# across all 8 sheets, per the scheduled price-refresh job.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1730.5555
$ws.Range("I80").Value = 674.4
$ws.Range("J80").Value = 3050.75
$ws.Range("K80").Value = 2023.2
$ws.Range("L80").Value = 9152.25
$ws.Range("M80").Value = -1025.2
$ws.Range("N80").Value = -11148.25
$ws.Range("H83").Value = 1730.5555
$ws.Range("I83").Value = 674.4
$ws.Range("J83").Value = 3050.75
$ws.Range("K83").Value = 6069.599999999999
$ws.Range("L83").Value = 27456.75
$ws.Range("M83").Value = -1077.599999999999
$ws.Range("N83").Value = -37440.75
$ws.Range("H98").Value = 533849.6
$ws.Range("I98").Value = 658249.5
$ws.Range("J98").Value = 5150
$ws.Range("K98").Value = 658249.5
$ws.Range("L98").Value = 5150
$ws.Range("M98").Value = -656751.5
$ws.Range("N98").Value = -8146
$ws.Range("H122").Value = 533849.6
$ws.Range("I122").Value = 658249.5
$ws.Range("J122").Value = 5150
$ws.Range("K122").Value = 1974748.5
$ws.Range("L122").Value = 15450
$ws.Range("M122").Value = -1972298.5
$ws.Range("N122").Value = -20350
$ws.Range("H135").Value = 1799.1052
$ws.Range("I135").Value = 1981.9286
$ws.Range("J135").Value = 1287.2
$ws.Range("K135").Value = 17837.3574
$ws.Range("L135").Value = 11584.8
$ws.Range("M135").Value = -15302.3574
$ws.Range("N135").Value = -16654.8
$ws.Range("H138").Value = 5378391.5
$ws.Range("J138").Value = 8622871
$ws.Range("L138").Value = 25868613
$ws.Range("N138").Value = -25878893

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 6950
$ws.Range("J43").Value = 6950
$ws.Range("L43").Value = 6950
$ws.Range("N43").Value = -7576
$ws.Range("H61").Value = 2358.7585
$ws.Range("I61").Value = 1607.0435
$ws.Range("K61").Value = 1607.0435
$ws.Range("M61").Value = -1395.0435
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H120").Value = 50000
$ws.Range("J120").Value = 50000
$ws.Range("L120").Value = 50000
$ws.Range("N120").Value = -59676
$ws.Range("H132").Value = 1823.4897
$ws.Range("I132").Value = 1395.1428
$ws.Range("J132").Value = 4393.5713
$ws.Range("K132").Value = 4185.428400000001
$ws.Range("L132").Value = 13180.7139
$ws.Range("M132").Value = -1655.428400000001
$ws.Range("N132").Value = -18240.7139
$ws.Range("H136").Value = 2358.7585
$ws.Range("I136").Value = 1607.0435
$ws.Range("K136").Value = 4821.1305
$ws.Range("M136").Value = -2271.1305

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1508.5416
$ws.Range("I20").Value = 1142.3636
$ws.Range("J20").Value = 1818.3846
$ws.Range("K20").Value = 1142.3636
$ws.Range("L20").Value = 1818.3846
$ws.Range("M20").Value = -895.3635999999999
$ws.Range("N20").Value = -2312.3846
$ws.Range("H105").Value = 10103445
$ws.Range("I105").Value = 11907121
$ws.Range("J105").Value = 2859.8
$ws.Range("K105").Value = 11907121
$ws.Range("L105").Value = 2859.8
$ws.Range("M105").Value = -11905374
$ws.Range("N105").Value = -6353.8
$ws.Range("H107").Value = 1328.8
$ws.Range("I107").Value = 1074.5
$ws.Range("K107").Value = 1074.5
$ws.Range("M107").Value = 845.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1458.75
$ws.Range("I31").Value = 1037.8379
$ws.Range("J31").Value = 1903.7142
$ws.Range("K31").Value = 1037.8379
$ws.Range("L31").Value = 1903.7142
$ws.Range("M31").Value = -742.8379
$ws.Range("N31").Value = -2493.7142
$ws.Range("H34").Value = 1458.75
$ws.Range("I34").Value = 1037.8379
$ws.Range("J34").Value = 1903.7142
$ws.Range("K34").Value = 1037.8379
$ws.Range("L34").Value = 1903.7142
$ws.Range("M34").Value = -835.8379
$ws.Range("N34").Value = -2307.7142
$ws.Range("H99").Value = 4168184
$ws.Range("I99").Value = 6251336
$ws.Range("J99").Value = 1879.8
$ws.Range("K99").Value = 6251336
$ws.Range("L99").Value = 1879.8
$ws.Range("M99").Value = -6249838
$ws.Range("N99").Value = -4875.8
$ws.Range("H126").Value = 4168184
$ws.Range("I126").Value = 6251336
$ws.Range("J126").Value = 1879.8
$ws.Range("K126").Value = 18754008
$ws.Range("L126").Value = 5639.4
$ws.Range("M126").Value = -18751538
$ws.Range("N126").Value = -10579.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1550.2222
$ws.Range("J34").Value = 2400
$ws.Range("L34").Value = 7200
$ws.Range("N34").Value = -7368
$ws.Range("H39").Value = 8290.736999999999
$ws.Range("J39").Value = 8290.736999999999
$ws.Range("L39").Value = 24872.211
$ws.Range("N39").Value = -25460.211
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H68").Value = 942.4400000000001
$ws.Range("I68").Value = 738.2769
$ws.Range("J68").Value = 1321.6
$ws.Range("K68").Value = 2214.8307
$ws.Range("L68").Value = 3964.8
$ws.Range("M68").Value = -1403.8307
$ws.Range("N68").Value = -5586.799999999999
$ws.Range("H71").Value = 942.4400000000001
$ws.Range("I71").Value = 738.2769
$ws.Range("J71").Value = 1321.6
$ws.Range("K71").Value = 6644.492099999999
$ws.Range("L71").Value = 11894.4
$ws.Range("M71").Value = -2588.492099999999
$ws.Range("N71").Value = -20006.4
$ws.Range("H132").Value = 1185.5
$ws.Range("I132").Value = 1038
$ws.Range("J132").Value = 1203.9375
$ws.Range("K132").Value = 9342
$ws.Range("L132").Value = 10835.4375
$ws.Range("M132").Value = -6812
$ws.Range("N132").Value = -15895.4375
$ws.Range("H133").Value = 9709.286
$ws.Range("I133").Value = 9160.833000000001
$ws.Range("K133").Value = 27482.499
$ws.Range("M133").Value = -22422.499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 20454.545
$ws.Range("J33").Value = 20454.545
$ws.Range("L33").Value = 20454.545
$ws.Range("N33").Value = -20958.545
$ws.Range("H51").Value = 57663
$ws.Range("J51").Value = 57663
$ws.Range("L51").Value = 57663
$ws.Range("N51").Value = -58681
$ws.Range("H92").Value = 13333.333
$ws.Range("J92").Value = 13333.333
$ws.Range("L92").Value = 13333.333
$ws.Range("N92").Value = -17077.333
$ws.Range("H97").Value = 740.9722
$ws.Range("I97").Value = 662.2222
$ws.Range("J97").Value = 819.7222
$ws.Range("K97").Value = 662.2222
$ws.Range("L97").Value = 819.7222
$ws.Range("M97").Value = -166.2222
$ws.Range("N97").Value = -1811.7222
$ws.Range("H109").Value = 71535
$ws.Range("J109").Value = 71535
$ws.Range("L109").Value = 71535
$ws.Range("N109").Value = -73615
$ws.Range("H122").Value = 794731.2
$ws.Range("J122").Value = 1800
$ws.Range("L122").Value = 5400
$ws.Range("N122").Value = -10300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3333.3333
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 3357.1428
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3357.1428
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -3581.1428
$ws.Range("H16").Value = 613.5
$ws.Range("I16").Value = 589.44446
$ws.Range("J16").Value = 685.6667
$ws.Range("K16").Value = 589.44446
$ws.Range("L16").Value = 685.6667
$ws.Range("M16").Value = -419.44446
$ws.Range("N16").Value = -1025.6667
$ws.Range("H105").Value = 50615
$ws.Range("J105").Value = 50615
$ws.Range("L105").Value = 50615
$ws.Range("N105").Value = -57603
$ws.Range("H126").Value = 3333.3333
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 3357.1428
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 10071.4284
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -15011.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H123").Value = 29745.6
$ws.Range("J123").Value = 29745.6
$ws.Range("L123").Value = 29745.6
$ws.Range("N123").Value = -39545.6
$ws.Range("H136").Value = 10449116
$ws.Range("I136").Value = 15197546
$ws.Range("J136").Value = 2569.1
$ws.Range("K136").Value = 45592638
$ws.Range("L136").Value = 7707.299999999999
$ws.Range("M136").Value = -45590088
$ws.Range("N136").Value = -12807.3

